# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.344.69"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "1.933.82"
$ws.Range("E3").Value = "  +1.04%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.53%  "
$ws.Range("E5").Value = "  +2.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7198"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3281"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "27.77"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07271"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8056"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08092"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.17%  "
$ws.Range("D13").Value = "1.934.56"
$ws.Range("E13").Value = "  +1.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.430"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.63"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.52%  "
$ws.Range("D17").Value = "30.339.28"
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008231"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.77%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "253.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.810"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.45%  "
$ws.Range("D21").Value = "2.187.50"
$ws.Range("E21").Value = "  +1.11%  "
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.947"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.721"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.40"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.344"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1293"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.360"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.545"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.442"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.202"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05232"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.270"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7493"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.761"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01966"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.61%  "
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "79.20"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.453"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4551"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.035"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8434"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.815"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.456"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.80"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4193"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.36%  "
$ws.Range("E51").Value = "  +2.67%  "
